# Add 2022-Q3 data
# -------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for the 2022-Q3 summary, shifting
#    the existing quarters down by one row, and renumber the index
#    column (A) for all the shifted rows.
# -------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(2).Insert()

# Carry the bold/bordered "index" cell formatting from A3 onto the new A2,
# and reset B2:D2 back to the plain/default formatting (Insert() copies the
# row-above's blank formatting into the new row by default).
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("B2:D2").Style = "Normal"

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.04

$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2
$ws1.Range("A5").Value = 3
$ws1.Range("A6").Value = 4
$ws1.Range("A7").Value = 5
$ws1.Range("A8").Value = 6
$ws1.Range("A9").Value = 7

# -------------------------------------------------------------------
# 2. Create the new "2022-Q3" worksheet (holdings detail) by copying the
#    "2022-Q2" sheet (so formatting/layout matches exactly) and placing
#    the copy right before it, then trim/update the data to the
#    2022-Q3 numbers.
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Copy($ws2)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Original copied rows: 2=168501, 3=009954, 4=001829, 5=002123, 6=004192, 7=004193
# Keep only the 009954 and 002123 rows (which already match the 2022-Q3 fund
# code/name), dropping the rest.
$newSheet.Rows.Item(2).Delete()
$newSheet.Rows.Item(3).Delete()
$newSheet.Rows.Item(4).Resize(2).Delete()

# Fix the numeric index/rank columns.
$newSheet.Range("A2").Value = 0
$newSheet.Range("H2").Value = 7
$newSheet.Range("A3").Value = 1
$newSheet.Range("H3").Value = 2

# Fix the text-valued metric columns (fund size / position / etc.), forcing
# them to stay text (matching the original inlineStr cells) instead of being
# auto-converted to numbers.
$newSheet.Range("D2:G3").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.55"
$newSheet.Range("E2").Value = "93.35"
$newSheet.Range("F2").Value = "5.68"
$newSheet.Range("G2").Value = "0.0312"
$newSheet.Range("D3").Value = "0.15"
$newSheet.Range("E3").Value = "88.55"
$newSheet.Range("F3").Value = "5.41"
$newSheet.Range("G3").Value = "0.0081"
$newSheet.Range("D2:G3").Style = "Normal"
